# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit: output generated at 7921097).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 270
    4  = 11162
    5  = 10399
    13 = 10430
    14 = 2221
    16 = 2444
    20 = 396
    21 = 11062
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
